$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSLXPPTF5L%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T090003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJGMEQCIByNgyG0wrR2KKpwVEVxp0qNKwd2eE1LZUo8AILnLcw9AiAojFkOOtGUz4MQ19qPA%2Fiat3117fI%2B7c8lVt25pNXizCrmAwiF%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM0frgMjiSYxgjpr7HKroDa7cdZ4E00FU%2FGVpS4encjtSDCUiIpLlhGUPqALd30WCtTpSOqOTrkSxzugUTh7SBcj1Lk8zBSD6KQ8HLxXrMJhKRDXkjrZCJBbs7fSIjDTNChpSfripixTeZ7t4D2LAN%2F9aAoVsUXx3ttPk4X9qTtahDL8uZlJYAtHn3XjGEaxUvV5liPV76Zt5g5cavEUkjHpF4ANgWEbmBAR8J45WxZg6cEYjJiiUxFATYGKXR%2FqP9RFNxSOOj%2FWvwSSyg6rjSw1bJGTjr%2FwkDjNTaRthUfQpRrGUm4a%2B6ezscUhjuTXTEaGD7OPR2j2bpFzznQP4w9NRIMHcFfDnWjro2SFN7FYWIi8VkaPc02wdf8Omv%2FsrZKVnNGLMbHSVKJDop6ERLBURHCvyeCNEJ7Kh0kA6dROAABJSVGILYTtZD74nlnDDNsJu3nx2iscpEm%2FfKwgwtP8e6M0YSFofa%2B201q%2BKf4D9MPiJ2PT0IWnRdFNyMtHYFufI4TIE3YOmo6ALsk9CtdgiUvWCtuiZ82hZmtudLYc6uv9gkYFi84D6A8DRXKU1hY9zzIKuJREYR5yFz7isct0%2FAg6a3RYMAOzC%2FiaO%2FBjqmATOYI%2FYR%2B9o46ZcKxgoMQPdKP8B9xS61DGgRun49TZ3lihG%2BKl8pyggLbpYcYW2gvGPDIVnCsbrtG%2FgDBQB7b60ItGso5LeAzOiTvvRc8y%2Fl667zX3oCmBNpLYhdQerjD4WJXwCRu%2BR6009VTYMEYQ9GYTHUePdfzsmvrZ8d%2B%2BUWS7nDZl%2F40%2Bm6ufWvB8l2yGAaSOZhB%2FYI3VZ6Z4sOrO1i5%2BEnpFI%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2017%3A00%3A03%20GMT&X-Amz-Signature=b9f2caedceaabd29f8818930900cfc40dda28a564290d945ed6ec6c281ae135e'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSLXPPTF5L%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T090003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJGMEQCIByNgyG0wrR2KKpwVEVxp0qNKwd2eE1LZUo8AILnLcw9AiAojFkOOtGUz4MQ19qPA%2Fiat3117fI%2B7c8lVt25pNXizCrmAwiF%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM0frgMjiSYxgjpr7HKroDa7cdZ4E00FU%2FGVpS4encjtSDCUiIpLlhGUPqALd30WCtTpSOqOTrkSxzugUTh7SBcj1Lk8zBSD6KQ8HLxXrMJhKRDXkjrZCJBbs7fSIjDTNChpSfripixTeZ7t4D2LAN%2F9aAoVsUXx3ttPk4X9qTtahDL8uZlJYAtHn3XjGEaxUvV5liPV76Zt5g5cavEUkjHpF4ANgWEbmBAR8J45WxZg6cEYjJiiUxFATYGKXR%2FqP9RFNxSOOj%2FWvwSSyg6rjSw1bJGTjr%2FwkDjNTaRthUfQpRrGUm4a%2B6ezscUhjuTXTEaGD7OPR2j2bpFzznQP4w9NRIMHcFfDnWjro2SFN7FYWIi8VkaPc02wdf8Omv%2FsrZKVnNGLMbHSVKJDop6ERLBURHCvyeCNEJ7Kh0kA6dROAABJSVGILYTtZD74nlnDDNsJu3nx2iscpEm%2FfKwgwtP8e6M0YSFofa%2B201q%2BKf4D9MPiJ2PT0IWnRdFNyMtHYFufI4TIE3YOmo6ALsk9CtdgiUvWCtuiZ82hZmtudLYc6uv9gkYFi84D6A8DRXKU1hY9zzIKuJREYR5yFz7isct0%2FAg6a3RYMAOzC%2FiaO%2FBjqmATOYI%2FYR%2B9o46ZcKxgoMQPdKP8B9xS61DGgRun49TZ3lihG%2BKl8pyggLbpYcYW2gvGPDIVnCsbrtG%2FgDBQB7b60ItGso5LeAzOiTvvRc8y%2Fl667zX3oCmBNpLYhdQerjD4WJXwCRu%2BR6009VTYMEYQ9GYTHUePdfzsmvrZ8d%2B%2BUWS7nDZl%2F40%2Bm6ufWvB8l2yGAaSOZhB%2FYI3VZ6Z4sOrO1i5%2BEnpFI%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2017%3A00%3A03%20GMT&X-Amz-Signature=7e22e5c009ca1564764defdeca9e545eba270ec8d66fea2edd17b2bc9ea74422'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSLXPPTF5L%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T090003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJGMEQCIByNgyG0wrR2KKpwVEVxp0qNKwd2eE1LZUo8AILnLcw9AiAojFkOOtGUz4MQ19qPA%2Fiat3117fI%2B7c8lVt25pNXizCrmAwiF%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM0frgMjiSYxgjpr7HKroDa7cdZ4E00FU%2FGVpS4encjtSDCUiIpLlhGUPqALd30WCtTpSOqOTrkSxzugUTh7SBcj1Lk8zBSD6KQ8HLxXrMJhKRDXkjrZCJBbs7fSIjDTNChpSfripixTeZ7t4D2LAN%2F9aAoVsUXx3ttPk4X9qTtahDL8uZlJYAtHn3XjGEaxUvV5liPV76Zt5g5cavEUkjHpF4ANgWEbmBAR8J45WxZg6cEYjJiiUxFATYGKXR%2FqP9RFNxSOOj%2FWvwSSyg6rjSw1bJGTjr%2FwkDjNTaRthUfQpRrGUm4a%2B6ezscUhjuTXTEaGD7OPR2j2bpFzznQP4w9NRIMHcFfDnWjro2SFN7FYWIi8VkaPc02wdf8Omv%2FsrZKVnNGLMbHSVKJDop6ERLBURHCvyeCNEJ7Kh0kA6dROAABJSVGILYTtZD74nlnDDNsJu3nx2iscpEm%2FfKwgwtP8e6M0YSFofa%2B201q%2BKf4D9MPiJ2PT0IWnRdFNyMtHYFufI4TIE3YOmo6ALsk9CtdgiUvWCtuiZ82hZmtudLYc6uv9gkYFi84D6A8DRXKU1hY9zzIKuJREYR5yFz7isct0%2FAg6a3RYMAOzC%2FiaO%2FBjqmATOYI%2FYR%2B9o46ZcKxgoMQPdKP8B9xS61DGgRun49TZ3lihG%2BKl8pyggLbpYcYW2gvGPDIVnCsbrtG%2FgDBQB7b60ItGso5LeAzOiTvvRc8y%2Fl667zX3oCmBNpLYhdQerjD4WJXwCRu%2BR6009VTYMEYQ9GYTHUePdfzsmvrZ8d%2B%2BUWS7nDZl%2F40%2Bm6ufWvB8l2yGAaSOZhB%2FYI3VZ6Z4sOrO1i5%2BEnpFI%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2017%3A00%3A03%20GMT&X-Amz-Signature=ed47db6649106c8676d33f05082075279afe575ef7830abe59c3e392bca6ca03'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSLXPPTF5L%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T090003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJGMEQCIByNgyG0wrR2KKpwVEVxp0qNKwd2eE1LZUo8AILnLcw9AiAojFkOOtGUz4MQ19qPA%2Fiat3117fI%2B7c8lVt25pNXizCrmAwiF%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM0frgMjiSYxgjpr7HKroDa7cdZ4E00FU%2FGVpS4encjtSDCUiIpLlhGUPqALd30WCtTpSOqOTrkSxzugUTh7SBcj1Lk8zBSD6KQ8HLxXrMJhKRDXkjrZCJBbs7fSIjDTNChpSfripixTeZ7t4D2LAN%2F9aAoVsUXx3ttPk4X9qTtahDL8uZlJYAtHn3XjGEaxUvV5liPV76Zt5g5cavEUkjHpF4ANgWEbmBAR8J45WxZg6cEYjJiiUxFATYGKXR%2FqP9RFNxSOOj%2FWvwSSyg6rjSw1bJGTjr%2FwkDjNTaRthUfQpRrGUm4a%2B6ezscUhjuTXTEaGD7OPR2j2bpFzznQP4w9NRIMHcFfDnWjro2SFN7FYWIi8VkaPc02wdf8Omv%2FsrZKVnNGLMbHSVKJDop6ERLBURHCvyeCNEJ7Kh0kA6dROAABJSVGILYTtZD74nlnDDNsJu3nx2iscpEm%2FfKwgwtP8e6M0YSFofa%2B201q%2BKf4D9MPiJ2PT0IWnRdFNyMtHYFufI4TIE3YOmo6ALsk9CtdgiUvWCtuiZ82hZmtudLYc6uv9gkYFi84D6A8DRXKU1hY9zzIKuJREYR5yFz7isct0%2FAg6a3RYMAOzC%2FiaO%2FBjqmATOYI%2FYR%2B9o46ZcKxgoMQPdKP8B9xS61DGgRun49TZ3lihG%2BKl8pyggLbpYcYW2gvGPDIVnCsbrtG%2FgDBQB7b60ItGso5LeAzOiTvvRc8y%2Fl667zX3oCmBNpLYhdQerjD4WJXwCRu%2BR6009VTYMEYQ9GYTHUePdfzsmvrZ8d%2B%2BUWS7nDZl%2F40%2Bm6ufWvB8l2yGAaSOZhB%2FYI3VZ6Z4sOrO1i5%2BEnpFI%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2017%3A00%3A03%20GMT&X-Amz-Signature=0b5d1b791e22ed800f457ef6f62a948663da906de2a869c5f64e2ba0436274d9'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSLXPPTF5L%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T090003Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBwaCXVzLXdlc3QtMiJGMEQCIByNgyG0wrR2KKpwVEVxp0qNKwd2eE1LZUo8AILnLcw9AiAojFkOOtGUz4MQ19qPA%2Fiat3117fI%2B7c8lVt25pNXizCrmAwiF%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM0frgMjiSYxgjpr7HKroDa7cdZ4E00FU%2FGVpS4encjtSDCUiIpLlhGUPqALd30WCtTpSOqOTrkSxzugUTh7SBcj1Lk8zBSD6KQ8HLxXrMJhKRDXkjrZCJBbs7fSIjDTNChpSfripixTeZ7t4D2LAN%2F9aAoVsUXx3ttPk4X9qTtahDL8uZlJYAtHn3XjGEaxUvV5liPV76Zt5g5cavEUkjHpF4ANgWEbmBAR8J45WxZg6cEYjJiiUxFATYGKXR%2FqP9RFNxSOOj%2FWvwSSyg6rjSw1bJGTjr%2FwkDjNTaRthUfQpRrGUm4a%2B6ezscUhjuTXTEaGD7OPR2j2bpFzznQP4w9NRIMHcFfDnWjro2SFN7FYWIi8VkaPc02wdf8Omv%2FsrZKVnNGLMbHSVKJDop6ERLBURHCvyeCNEJ7Kh0kA6dROAABJSVGILYTtZD74nlnDDNsJu3nx2iscpEm%2FfKwgwtP8e6M0YSFofa%2B201q%2BKf4D9MPiJ2PT0IWnRdFNyMtHYFufI4TIE3YOmo6ALsk9CtdgiUvWCtuiZ82hZmtudLYc6uv9gkYFi84D6A8DRXKU1hY9zzIKuJREYR5yFz7isct0%2FAg6a3RYMAOzC%2FiaO%2FBjqmATOYI%2FYR%2B9o46ZcKxgoMQPdKP8B9xS61DGgRun49TZ3lihG%2BKl8pyggLbpYcYW2gvGPDIVnCsbrtG%2FgDBQB7b60ItGso5LeAzOiTvvRc8y%2Fl667zX3oCmBNpLYhdQerjD4WJXwCRu%2BR6009VTYMEYQ9GYTHUePdfzsmvrZ8d%2B%2BUWS7nDZl%2F40%2Bm6ufWvB8l2yGAaSOZhB%2FYI3VZ6Z4sOrO1i5%2BEnpFI%3D&X-Amz-SignedHeaders=host&response-expires=Sun%2C%2030%20Mar%202025%2017%3A00%3A03%20GMT&X-Amz-Signature=cf5a310251a948a85535ff39107d86e59182da0855d84c6c7a19bcfd81bdfd23'
